# "TC for User reglogin added"
#
# The Login_TestCases sheet had a block of stale, empty leftover rows
# (rows 2-20) sitting above the real table (header in row 21, data rows
# 22-33). Clean those out so the table starts right at the top, then add
# the missing Test Steps / Precondition text for the "Verify login button
# functionality" test case (TC_LOGIN_08), which previously had no data in
# those two columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login_TestCases")
$ws.Activate()

# Remove the 19 blank rows above the header row, shifting the table
# (old row 21 header / rows 22-33 data) up to rows 2/3-14.
$ws.Rows("2:20").Delete()

# The table header is now row 2, and TC_LOGIN_08 ("Verify login button
# functionality") is now row 10. Fill in its Steps (column G) and
# Precondition (column F) - previously blank.
$ws.Range("G10").Value = "1.Enter Email.`n2.Enter Password`n3.Clik login "
$ws.Range("F10").Value = "Home Page to be displayed"

# Match the wrapped-text styling used by the Steps column elsewhere in the
# table, and size the row to fit the new 3-line entry.
$ws.Range("G10").WrapText = $true
$ws.Rows(10).RowHeight = 45

# Leave the view focused on the (now relocated) header row.
$ws.Rows(2).Select()
